# Horarios actualizados Linea 141 - 726
# Refresh the scraped-schedule data on all three sheets:
#   LP1912      (sheet 1)
#   LP1912-215  (sheet 2)
#   6203-6173   (sheet 3)
# "Ultima actualizacion" / "Total filas" headers + row data are updated to the
# new scrape taken at 04:44:38, including newly appended rows.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# ===================== Sheet 1: LP1912 =====================
$ws1.Range("A2").Value = "Última actualización: 04:44:38"
$ws1.Range("A3").Value = "Total filas: 32"

# Rows 6-18 are unchanged from the previous scrape.
# Rows 19-25, 27-37 are updated/new (row 26 stays the same).
$ws1.Cells.Item(19,1).Value = "04:44:38"
$ws1.Cells.Item(19,2).Value = "04:45"
$ws1.Cells.Item(19,3).Value = "215A_EL PATO"
$ws1.Cells.Item(19,4).Value = 1
$ws1.Cells.Item(19,5).Value = "LP1912"

$ws1.Cells.Item(20,1).Value = "04:44:38"
$ws1.Cells.Item(20,2).Value = "04:53"
$ws1.Cells.Item(20,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(20,4).Value = 9
$ws1.Cells.Item(20,5).Value = "LP1912"

$ws1.Cells.Item(21,1).Value = "04:44:38"
$ws1.Cells.Item(21,2).Value = "05:16"
$ws1.Cells.Item(21,3).Value = "17_ROMERO"
$ws1.Cells.Item(21,4).Value = 32
$ws1.Cells.Item(21,5).Value = "LP1912"

$ws1.Cells.Item(22,1).Value = "04:44:38"
$ws1.Cells.Item(22,2).Value = "05:22"
$ws1.Cells.Item(22,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(22,4).Value = 38
$ws1.Cells.Item(22,5).Value = "LP1912"

$ws1.Cells.Item(23,1).Value = "04:44:38"
$ws1.Cells.Item(23,2).Value = "05:34"
$ws1.Cells.Item(23,3).Value = "215B_EL PATO"
$ws1.Cells.Item(23,4).Value = 50
$ws1.Cells.Item(23,5).Value = "LP1912"

$ws1.Cells.Item(24,1).Value = "04:17:03"
$ws1.Cells.Item(24,2).Value = "05:35"
$ws1.Cells.Item(24,3).Value = "215B_EL PATO"
$ws1.Cells.Item(24,4).Value = 78
$ws1.Cells.Item(24,5).Value = "LP1912"

$ws1.Cells.Item(25,1).Value = "03:42:43"
$ws1.Cells.Item(25,2).Value = "05:35"
$ws1.Cells.Item(25,3).Value = "14_ABASTO"
$ws1.Cells.Item(25,4).Value = 113
$ws1.Cells.Item(25,5).Value = "LP1912"

# Row 26 unchanged: 04:17:03 | 05:36 | 14_ABASTO | 79 | LP1912

$ws1.Cells.Item(27,1).Value = "04:44:38"
$ws1.Cells.Item(27,2).Value = "05:46"
$ws1.Cells.Item(27,3).Value = "15_ABASTO"
$ws1.Cells.Item(27,4).Value = 62
$ws1.Cells.Item(27,5).Value = "LP1912"

$ws1.Cells.Item(28,1).Value = "04:44:38"
$ws1.Cells.Item(28,2).Value = "05:54"
$ws1.Cells.Item(28,3).Value = "10_OLMOS"
$ws1.Cells.Item(28,4).Value = 70
$ws1.Cells.Item(28,5).Value = "LP1912"

$ws1.Cells.Item(29,1).Value = "04:44:38"
$ws1.Cells.Item(29,2).Value = "06:04"
$ws1.Cells.Item(29,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(29,4).Value = 80
$ws1.Cells.Item(29,5).Value = "LP1912"

$ws1.Cells.Item(30,1).Value = "04:17:03"
$ws1.Cells.Item(30,2).Value = "06:05"
$ws1.Cells.Item(30,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(30,4).Value = 108
$ws1.Cells.Item(30,5).Value = "LP1912"

$ws1.Cells.Item(31,1).Value = "04:44:38"
$ws1.Cells.Item(31,2).Value = "06:11"
$ws1.Cells.Item(31,3).Value = "215A_EL PATO"
$ws1.Cells.Item(31,4).Value = 87
$ws1.Cells.Item(31,5).Value = "LP1912"

$ws1.Cells.Item(32,1).Value = "04:17:03"
$ws1.Cells.Item(32,2).Value = "06:12"
$ws1.Cells.Item(32,3).Value = "215A_EL PATO"
$ws1.Cells.Item(32,4).Value = 115
$ws1.Cells.Item(32,5).Value = "LP1912"

$ws1.Cells.Item(33,1).Value = "04:44:38"
$ws1.Cells.Item(33,2).Value = "06:14"
$ws1.Cells.Item(33,3).Value = "225_HARAS DEL SUR"
$ws1.Cells.Item(33,4).Value = 90
$ws1.Cells.Item(33,5).Value = "LP1912"

$ws1.Cells.Item(34,1).Value = "04:44:38"
$ws1.Cells.Item(34,2).Value = "06:21"
$ws1.Cells.Item(34,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(34,4).Value = 97
$ws1.Cells.Item(34,5).Value = "LP1912"

$ws1.Cells.Item(35,1).Value = "04:44:38"
$ws1.Cells.Item(35,2).Value = "06:27"
$ws1.Cells.Item(35,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(35,4).Value = 103
$ws1.Cells.Item(35,5).Value = "LP1912"

$ws1.Cells.Item(36,1).Value = "04:44:38"
$ws1.Cells.Item(36,2).Value = "06:29"
$ws1.Cells.Item(36,3).Value = "86_EST CHICA-ESC AGRARIA"
$ws1.Cells.Item(36,4).Value = 105
$ws1.Cells.Item(36,5).Value = "LP1912"

$ws1.Cells.Item(37,1).Value = "04:44:38"
$ws1.Cells.Item(37,2).Value = "06:31"
$ws1.Cells.Item(37,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(37,4).Value = 107
$ws1.Cells.Item(37,5).Value = "LP1912"

# ===================== Sheet 2: LP1912-215 =====================
$ws2.Range("A2").Value = "Última actualización: 04:44:38"
$ws2.Range("A3").Value = "Total filas: 13"

# Rows 6-13 unchanged from the previous scrape.
$ws2.Cells.Item(14,1).Value = "04:44:38"
$ws2.Cells.Item(14,2).Value = "04:45"
$ws2.Cells.Item(14,3).Value = "215A_EL PATO"
$ws2.Cells.Item(14,4).Value = 1
$ws2.Cells.Item(14,5).Value = "LP1912"

$ws2.Cells.Item(15,1).Value = "04:44:38"
$ws2.Cells.Item(15,2).Value = "05:34"
$ws2.Cells.Item(15,3).Value = "215B_EL PATO"
$ws2.Cells.Item(15,4).Value = 50
$ws2.Cells.Item(15,5).Value = "LP1912"

# Row 16 unchanged: 04:17:03 | 05:35 | 215B_EL PATO | 78 | LP1912

$ws2.Cells.Item(17,1).Value = "04:44:38"
$ws2.Cells.Item(17,2).Value = "06:11"
$ws2.Cells.Item(17,3).Value = "215A_EL PATO"
$ws2.Cells.Item(17,4).Value = 87
$ws2.Cells.Item(17,5).Value = "LP1912"

$ws2.Cells.Item(18,1).Value = "04:17:03"
$ws2.Cells.Item(18,2).Value = "06:12"
$ws2.Cells.Item(18,3).Value = "215A_EL PATO"
$ws2.Cells.Item(18,4).Value = 115
$ws2.Cells.Item(18,5).Value = "LP1912"

# ===================== Sheet 3: 6203-6173 =====================
$ws3.Range("A2").Value = "Última actualización: 04:44:38"
$ws3.Range("A3").Value = "Total filas: 5"

# Row 6 unchanged: 00:05:23 | 00:08 | 215A_LA PLATA | 3 | L6173
$ws3.Cells.Item(7,1).Value = "04:44:38"
$ws3.Cells.Item(7,2).Value = "05:43"
$ws3.Cells.Item(7,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(7,4).Value = 59
$ws3.Cells.Item(7,5).Value = "L6173"

$ws3.Cells.Item(8,1).Value = "04:17:03"
$ws3.Cells.Item(8,2).Value = "05:44"
$ws3.Cells.Item(8,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(8,4).Value = 87
$ws3.Cells.Item(8,5).Value = "L6173"

$ws3.Cells.Item(9,1).Value = "04:44:38"
$ws3.Cells.Item(9,2).Value = "06:08"
$ws3.Cells.Item(9,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(9,4).Value = 84
$ws3.Cells.Item(9,5).Value = "L6173"

$ws3.Cells.Item(10,1).Value = "04:44:38"
$ws3.Cells.Item(10,2).Value = "06:32"
$ws3.Cells.Item(10,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(10,4).Value = 108
$ws3.Cells.Item(10,5).Value = "L6203"
